$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to be treated as literal text so Excel doesn't
    # reinterpret numeric-looking strings (e.g. "1.00", "61.010.17")
    # as numbers and strip formatting / introduce float rounding.
    $range.NumberFormat = "@"
    $range.Value = $value
    # Restore the default "Normal" style so no stray number-format
    # style gets left behind on cells that didn't have one before.
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "61.010.17"
Set-TextValue $ws.Range("E2") "  +0.44%  "
Set-TextValue $ws.Range("D3") "2.915.47"
Set-TextValue $ws.Range("E3") "  +0.47%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  -0.08%  "
Set-TextValue $ws.Range("D5") "590.28"
Set-TextValue $ws.Range("E5") "  +0.72%  "
Set-TextValue $ws.Range("D6") "144.97"
Set-TextValue $ws.Range("E6") "  +0.43%  "
Set-TextValue $ws.Range("E7") "  -0.03%  "
Set-TextValue $ws.Range("E8") "  +0.45%  "
Set-TextValue $ws.Range("E9") "  +3.87%  "
Set-TextValue $ws.Range("E10") "  -2.03%  "
Set-TextValue $ws.Range("E11") "  -1.48%  "
Set-TextValue $ws.Range("E12") "  -0.34%  "
Set-TextValue $ws.Range("D13") "33.39"
Set-TextValue $ws.Range("E13") "  +0.07%  "
Set-TextValue $ws.Range("E14") "  -0.30%  "
Set-TextValue $ws.Range("D15") "3.397.42"
Set-TextValue $ws.Range("E15") "  +0.25%  "
Set-TextValue $ws.Range("D16") "60.883.57"
Set-TextValue $ws.Range("E16") "  +0.24%  "
Set-TextValue $ws.Range("E17") "  -0.33%  "
Set-TextValue $ws.Range("D18") "2.915.67"
Set-TextValue $ws.Range("E18") "  +0.31%  "
Set-TextValue $ws.Range("D19") "433.15"
Set-TextValue $ws.Range("E19") "  +1.25%  "
Set-TextValue $ws.Range("D20") "13.36"
Set-TextValue $ws.Range("E20") "  -1.23%  "
Set-TextValue $ws.Range("E21") "  -0.65%  "
Set-TextValue $ws.Range("D22") "7.11"
Set-TextValue $ws.Range("E22") "  +0.24%  "
Set-TextValue $ws.Range("D23") "81.42"
Set-TextValue $ws.Range("E23") "  +1.10%  "
Set-TextValue $ws.Range("D24") "10.83"
Set-TextValue $ws.Range("E25") "  -0.69%  "
Set-TextValue $ws.Range("D26") "11.78"
Set-TextValue $ws.Range("E26") "  -0.91%  "
Set-TextValue $ws.Range("E27") "  -0.01%  "
Set-TextValue $ws.Range("E28") "  +4.75%  "
Set-TextValue $ws.Range("E29") "  -0.68%  "
Set-TextValue $ws.Range("D30") "6.96"
Set-TextValue $ws.Range("E30") "  -3.56%  "
Set-TextValue $ws.Range("D31") "26.47"
Set-TextValue $ws.Range("E31") "  +0.29%  "
Set-TextValue $ws.Range("E32") "  +2.91%  "
Set-TextValue $ws.Range("E33") "  -0.06%  "
Set-TextValue $ws.Range("D34") "0.0₃0869"
Set-TextValue $ws.Range("E34") "  -0.64%  "
Set-TextValue $ws.Range("E35") "  +0.45%  "
Set-TextValue $ws.Range("E36") "  +0.24%  "
Set-TextValue $ws.Range("D37") "3.01"
Set-TextValue $ws.Range("E37") "  -0.09%  "
Set-TextValue $ws.Range("E38") "  -1.04%  "
Set-TextValue $ws.Range("E39") "  -3.68%  "
Set-TextValue $ws.Range("E40") "  -0.57%  "
Set-TextValue $ws.Range("B41") "Arweave"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue $ws.Range("D41") "41.89"
Set-TextValue $ws.Range("E41") "  +1.40%  "
Set-TextValue $ws.Range("B42") "TheGraph"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D42") "0.288"
Set-TextValue $ws.Range("E42") "  -2.61%  "
Set-TextValue $ws.Range("D43") "375.42"
Set-TextValue $ws.Range("D44") "0.0346"
Set-TextValue $ws.Range("E44") "  -1.01%  "
Set-TextValue $ws.Range("D45") "2.691.38"
Set-TextValue $ws.Range("D46") "133.03"
Set-TextValue $ws.Range("E46") "  +0.41%  "
Set-TextValue $ws.Range("D48") "23.77"
Set-TextValue $ws.Range("E48") "  -2.11%  "
Set-TextValue $ws.Range("D49") "0.106"
Set-TextValue $ws.Range("E49") "  -0.46%  "
Set-TextValue $ws.Range("E50") "  -2.43%  "
Set-TextValue $ws.Range("D51") "0.123"
Set-TextValue $ws.Range("E51") "  -0.40%  "
